# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap the country labels whose case counts crossed each other ---
# Sweden (A23) overtook Ireland (A22)
$ws.Range("A22").Value = "Suecia"
$ws.Range("A23").Value = "Irlanda"

# Azerbaijan (A73) overtook Armenia (A72)
$ws.Range("A72").Value = "Azerbaiyan"
$ws.Range("A73").Value = "Armenia"

# --- Update the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 14:22"

# --- Row 4: Estados Unidos ---
$ws.Range("E4").Value = 717358
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 47684

# --- Row 8: Alemania ---
$ws.Range("E8").Value = 42154
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 5319

# --- Row 17: Paises Bajos ---
$ws.Range("B17").Value = 35729
$ws.Range("C17").Value = 887
$ws.Range("E17").Value = 31302
$ws.Range("G17").Value = 123
$ws.Range("H17").Value = 4177

# --- Row 19: Portugal ---
$ws.Range("B19").Value = 22353
$ws.Range("C19").Value = 371
$ws.Range("D19").Value = 1201
$ws.Range("E19").Value = 20332
$ws.Range("F19").Value = 204
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = 820

# --- Row 22: now Suecia (updated stats) ---
$ws.Range("B22").Value = 16755
$ws.Range("C22").Value = 751
$ws.Range("D22").Value = 550
$ws.Range("E22").Value = 14184
$ws.Range("F22").Value = 515
$ws.Range("G22").Value = 84
$ws.Range("H22").Value = 2021

# --- Row 23: now Irlanda (stats carried over) ---
$ws.Range("B23").Value = 16671
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 9233
$ws.Range("E23").Value = 6669
$ws.Range("F23").Value = 147
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 769

# --- Row 39: Catar ---
$ws.Range("F39").Value = 72

# --- Row 64: Kazajistan ---
$ws.Range("B64").Value = 2251
$ws.Range("C64").Value = 116
$ws.Range("E64").Value = 1695

# --- Row 66: Croacia ---
$ws.Range("B66").Value = 1981
$ws.Range("C66").Value = 31
$ws.Range("D66").Value = 883
$ws.Range("E66").Value = 1048
$ws.Range("G66").Value = 2
$ws.Range("H66").Value = 50

# --- Row 72: now Azerbaiyan (updated stats) ---
$ws.Range("B72").Value = 1548
$ws.Range("C72").Value = 30
$ws.Range("D72").Value = 948
$ws.Range("E72").Value = 580
$ws.Range("F72").Value = 14
$ws.Range("H72").Value = 20

# --- Row 73: now Armenia (stats carried over) ---
$ws.Range("B73").Value = 1523
$ws.Range("C73").Value = 50
$ws.Range("D73").Value = 659
$ws.Range("E73").Value = 840
$ws.Range("F73").Value = 10
$ws.Range("H73").Value = 24
